$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# Append a second "use case" block (rows 9-16) below the existing one
# (rows 1-8), mirroring its layout/styling exactly, for the new
# "View profile and statistics" use case.
# ------------------------------------------------------------------

# Row 10 mirrors row 2's merged title-bar cell (B2:C2) - merge first
# (while still blank/default-styled) then paste the formatting so the
# merge does not disturb the copied border style.
$ws.Range("B10:C10").Merge() | Out-Null
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B10:C10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value = "UC : View profile and statistics"

# Row 11 mirrors row 3 (Actor / System header cells).
$ws.Range("B3:C3").Copy() | Out-Null
$ws.Range("B11:C11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "Actor : User"
$ws.Range("C11").Value = "System : Card Game System"

# Row 12 mirrors row 4 (blank left cell / step description on right).
$ws.Range("B4:C4").Copy() | Out-Null
$ws.Range("B12:C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = "0.  System displays the main screen."

# Row 13 mirrors row 5 (two-step row).
$ws.Range("B5:C5").Copy() | Out-Null
$ws.Range("B13:C13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = "1.  TUCBW the user clicking on the profile button"
$ws.Range("C13").Value = "2. The user's profile is displayed with basic information."

# Row 14 mirrors row 6 (two-step row).
$ws.Range("B6:C6").Copy() | Out-Null
$ws.Range("B14:C14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = "3.  The user clicks the statistics button "
$ws.Range("C14").Value = "4.  Total statistics such as total chip count and total wins/losses are displayed, along with a list of games."

# Row 15 mirrors row 7 (two-step row).
$ws.Range("B7:C7").Copy() | Out-Null
$ws.Range("B15:C15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").Value = "5.  The user clicks on one of the individual games."
$ws.Range("C15").Value = "6.  Detailed statistics for that game are displayed, such as chips won or lost, wins and losses, and sessions played."

# Row 16 mirrors row 8 (closing step / bottom border row).
$ws.Range("B8:C8").Copy() | Out-Null
$ws.Range("B16:C16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = "7.  TUCEW the user reading the statistics and exploring further to additional game types"

$excel.CutCopyMode = 0

# Row heights (row 9 is the blank spacer row, matching row 1's).
$ws.Rows.Item(9).RowHeight = 15.75
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 45.75

$ws.Range("G15").Select()
